# Save the excel report:
#  - add a border-box around the new G1 header cell
#  - give the "Diff." column + the Optimizer sub-table a full thin border
#    (instead of the former left/right-only border)
#  - append a merged note row (A9:G11) explaining the seed investigation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: extend the bordered header box to the new G1 cell ---
$ws.Range("G1").Borders.LineStyle = 1

# --- Give the "Diff." header + values a full (4-sided) thin border ---
$ws.Range("G2:G8").Borders.LineStyle = 1

# --- Give the Optimizer sub-block (rows 7:8) a full thin border ---
$ws.Range("B7:D8").Borders.LineStyle = 1
$ws.Range("E7:F8").Borders.LineStyle = 1

# --- New note, merged across A9:G11 ---
$note = "Independent layer 0 and incremental layer 0 will have the same accuracy " + [char]10 + "if they have the same initializations"
$ws.Range("A9").Value = $note

$noteRange = $ws.Range("A9:G11")
$noteRange.Merge()
$noteRange.Borders.LineStyle = 1
$noteRange.HorizontalAlignment = -4108
$noteRange.VerticalAlignment = -4108

$ws.Range("A9").WrapText = $true

$wb.Save()
